$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Push the formatting of the current blank row (row 3, style index 1 /
#        Arial) down into the new blank row 4 before we touch anything else,
#        using a formats-only paste so the existing cell-format record is
#        reused on row 4, rather than synthesising a brand-new (but merely
#        equivalent) font/style entry.
$ws.Range("A3:D3").Copy()
$ws.Range("A4:D4").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- 2. Update row 2: "LUCAS CORTES" / "221" becomes "LUCAS" / "lucas"
$ws.Range("A2").Value = "LUCAS"
$ws.Range("B2").Value = "lucas"

# --- 3. Row 3 becomes a brand-new data record: LUCAS CORTES / 221 / 0 / 0,
#        with the default (no explicit) style.
$ws.Range("A3:D3").Style = "Normal"
$ws.Range("A3").Value = "LUCAS CORTES"
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "221"
$ws.Range("A3:D3").Style = "Normal"
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0

# --- 4. Row 4 is now the blank placeholder row (formatting already copied
#        above); make sure its cells are numeric/empty, matching the
#        original blank-row shape.
$ws.Range("A4:D4").ClearContents()
